$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CompleteBookDBTill4thMarch2022")

# Columns actually populated on existing data rows (T, U, W are intentionally
# left blank throughout the sheet, so we skip them for the new rows too).
$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","V","X","Y","Z","AA","AB","AC","AD","AE","AF")

# Copy the formatting (style) of the last existing data row (186) onto each
# of the three new rows, cell by cell, so the new rows keep the same look
# the rest of the table uses without inventing blank cells in unused columns.
$templateRow = 186
$newRowNumbers = @(187, 188, 189)
foreach ($newRow in $newRowNumbers) {
    foreach ($col in $cols) {
        $srcRange = $ws.Range($col + $templateRow)
        $dstRange = $ws.Range($col + $newRow)
        $srcRange.Copy()
        $dstRange.PasteSpecial(-4122)
    }
}

# Row 187 - "I Have A Dream"
$ws.Cells.Item(187, 1).Value = 186
$ws.Cells.Item(187, 2).Value = "Book"
$ws.Cells.Item(187, 3).Value = "I Have A Dream"
$ws.Cells.Item(187, 4).Value = "Yes"
$ws.Cells.Item(187, 5).Value = "Rashmi Bansal"
$ws.Cells.Item(187, 6).Value = "Stories of 20 Social Entrepruners"
$ws.Cells.Item(187, 7).Value = "New Arrivals"
$ws.Cells.Item(187, 8).Value = 349
$ws.Cells.Item(187, 9).Value = "NA"
$ws.Cells.Item(187, 10).Value = "Rs. 200"
$ws.Cells.Item(187, 11).Value = 1
$ws.Cells.Item(187, 12).Value = "1-1-E"
$ws.Cells.Item(187, 13).Value = "Paperback"
$ws.Cells.Item(187, 14).Value = 2011
$ws.Cells.Item(187, 15).Value = "English"
$ws.Cells.Item(187, 16).Value = "Self Help"
$ws.Cells.Item(187, 17).Value = "Entreprunership"
$ws.Cells.Item(187, 18).Value = "No"
$ws.Cells.Item(187, 19).Value = "NF"
$ws.Cells.Item(187, 22).Value = 0
$ws.Cells.Item(187, 24).Value = "Female"
$ws.Cells.Item(187, 25).Value = "978-93-80658-38-4"
$ws.Cells.Item(187, 26).Value = 1
$ws.Cells.Item(187, 27).Value = "First Floor"
$ws.Cells.Item(187, 28).Value = "Practical"
$ws.Cells.Item(187, 29).Value = "Yes"
$ws.Cells.Item(187, 30).Value = "Yes"
$ws.Cells.Item(187, 31).Value = 8.2
$ws.Cells.Item(187, 32).Value = "Good"

# Row 188 - "An Era Of Darkness"
$ws.Cells.Item(188, 1).Value = 187
$ws.Cells.Item(188, 2).Value = "Book"
$ws.Cells.Item(188, 3).Value = "An Era Of Darkness"
$ws.Cells.Item(188, 4).Value = "No"
$ws.Cells.Item(188, 5).Value = "Sashi Tharoor"
$ws.Cells.Item(188, 6).Value = "The British Empire In India"
$ws.Cells.Item(188, 7).Value = "New Arrivals"
$ws.Cells.Item(188, 8).Value = 333
$ws.Cells.Item(188, 9).Value = "Business Standard"
$ws.Cells.Item(188, 10).Value = "Rs. 699"
$ws.Cells.Item(188, 11).Value = 1
$ws.Cells.Item(188, 12).Value = "1-1-E"
$ws.Cells.Item(188, 13).Value = "Paperback"
$ws.Cells.Item(188, 14).Value = 2016
$ws.Cells.Item(188, 15).Value = "English"
$ws.Cells.Item(188, 16).Value = "History"
$ws.Cells.Item(188, 17).Value = "Indian History"
$ws.Cells.Item(188, 18).Value = "No"
$ws.Cells.Item(188, 19).Value = "NF"
$ws.Cells.Item(188, 22).Value = 0
$ws.Cells.Item(188, 24).Value = "Male"
$ws.Cells.Item(188, 25).Value = "978-93-83064-65-6"
$ws.Cells.Item(188, 26).Value = 1
$ws.Cells.Item(188, 27).Value = "First Floor"
$ws.Cells.Item(188, 28).Value = "Practical"
$ws.Cells.Item(188, 29).Value = "No"
$ws.Cells.Item(188, 30).Value = "No"
$ws.Cells.Item(188, 31).Value = 8.1
$ws.Cells.Item(188, 32).Value = "New"

# Row 189 - "Life Lessons For Loving The Way You Live"
$ws.Cells.Item(189, 1).Value = 188
$ws.Cells.Item(189, 2).Value = "Book"
$ws.Cells.Item(189, 3).Value = "Life Lessons For Loving The Way You Live"
$ws.Cells.Item(189, 4).Value = "Yes"
$ws.Cells.Item(189, 5).Value = "Jennifer Read Hawthorne"
$ws.Cells.Item(189, 6).Value = "7 essential ingredients for finding balance in life"
$ws.Cells.Item(189, 7).Value = "New Arrivals"
$ws.Cells.Item(189, 8).Value = 331
$ws.Cells.Item(189, 9).Value = "NA"
$ws.Cells.Item(189, 10).Value = "Rs. 2032"
$ws.Cells.Item(189, 11).Value = 1
$ws.Cells.Item(189, 12).Value = "1-1-E"
$ws.Cells.Item(189, 13).Value = "Paperback"
$ws.Cells.Item(189, 14).Value = 2008
$ws.Cells.Item(189, 15).Value = "English"
$ws.Cells.Item(189, 16).Value = "Self Help"
$ws.Cells.Item(189, 17).Value = "Self Love"
$ws.Cells.Item(189, 18).Value = "No"
$ws.Cells.Item(189, 19).Value = "NF"
$ws.Cells.Item(189, 22).Value = 0
$ws.Cells.Item(189, 24).Value = "Female"
$ws.Cells.Item(189, 25).Value = "978-81-89975-34-0"
$ws.Cells.Item(189, 26).Value = 1
$ws.Cells.Item(189, 27).Value = "First Floor"
$ws.Cells.Item(189, 28).Value = "Practical"
$ws.Cells.Item(189, 29).Value = "Yes"
$ws.Cells.Item(189, 30).Value = "Yes"
$ws.Cells.Item(189, 31).Value = 8.3
$ws.Cells.Item(189, 32).Value = "Good"
